$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell is set to Text format before assigning its value so that
# numeric-looking strings (with trailing zeros), percentages, dates, and
# comma-grouped numbers are preserved exactly as text, matching the source data.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '296.06'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '3.09%'
$ws.Range('F2').NumberFormat = '@'
$ws.Range('F2').Value = '14-1-2023'
$ws.Range('G2').NumberFormat = '@'
$ws.Range('G2').Value = '1'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '31.93'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '8.85%'
$ws.Range('F3').NumberFormat = '@'
$ws.Range('F3').Value = '14-1-2023'
$ws.Range('G3').NumberFormat = '@'
$ws.Range('G3').Value = '1'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.180'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '1.51%'
$ws.Range('F4').NumberFormat = '@'
$ws.Range('F4').Value = '14-1-2023'
$ws.Range('G4').NumberFormat = '@'
$ws.Range('G4').Value = '1'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.07357'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '11.05%'
$ws.Range('F5').NumberFormat = '@'
$ws.Range('F5').Value = '14-1-2023'
$ws.Range('G5').NumberFormat = '@'
$ws.Range('G5').Value = '1'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '7.666'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '4.54%'
$ws.Range('F6').NumberFormat = '@'
$ws.Range('F6').Value = '14-1-2023'
$ws.Range('G6').NumberFormat = '@'
$ws.Range('G6').Value = '1'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.648'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '6.87%'
$ws.Range('F7').NumberFormat = '@'
$ws.Range('F7').Value = '14-1-2023'
$ws.Range('G7').NumberFormat = '@'
$ws.Range('G7').Value = '1'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.439'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '5.35%'
$ws.Range('F8').NumberFormat = '@'
$ws.Range('F8').Value = '14-1-2023'
$ws.Range('G8').NumberFormat = '@'
$ws.Range('G8').Value = '1'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9031'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-2.59%'
$ws.Range('F9').NumberFormat = '@'
$ws.Range('F9').Value = '14-1-2023'
$ws.Range('G9').NumberFormat = '@'
$ws.Range('G9').Value = '1'
$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = 'One'
$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.01525'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '2,259.28%'
$ws.Range('F10').NumberFormat = '@'
$ws.Range('F10').Value = '14-1-2023'
$ws.Range('G10').NumberFormat = '@'
$ws.Range('G10').Value = '1'
$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1640'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '4.38%'
$ws.Range('F11').NumberFormat = '@'
$ws.Range('F11').Value = '14-1-2023'
$ws.Range('G11').NumberFormat = '@'
$ws.Range('G11').Value = '1'
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08132'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '24.39%'
$ws.Range('F12').NumberFormat = '@'
$ws.Range('F12').Value = '14-1-2023'
$ws.Range('G12').NumberFormat = '@'
$ws.Range('G12').Value = '1'
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'MandalaExchangeToken'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07736'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '1.78%'
$ws.Range('F13').NumberFormat = '@'
$ws.Range('F13').Value = '14-1-2023'
$ws.Range('G13').NumberFormat = '@'
$ws.Range('G13').Value = '1'
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.02987'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '2.62%'
$ws.Range('F14').NumberFormat = '@'
$ws.Range('F14').Value = '14-1-2023'
$ws.Range('G14').NumberFormat = '@'
$ws.Range('G14').Value = '1'
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.09336'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '3.83%'
$ws.Range('F15').NumberFormat = '@'
$ws.Range('F15').Value = '14-1-2023'
$ws.Range('G15').NumberFormat = '@'
$ws.Range('G15').Value = '1'
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.001519'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-4.81%'
$ws.Range('F16').NumberFormat = '@'
$ws.Range('F16').Value = '14-1-2023'
$ws.Range('G16').NumberFormat = '@'
$ws.Range('G16').Value = '1'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.006066'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-3.40%'
$ws.Range('F17').NumberFormat = '@'
$ws.Range('F17').Value = '14-1-2023'
$ws.Range('G17').NumberFormat = '@'
$ws.Range('G17').Value = '1'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.528'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '2.43%'
$ws.Range('F18').NumberFormat = '@'
$ws.Range('F18').Value = '14-1-2023'
$ws.Range('G18').NumberFormat = '@'
$ws.Range('G18').Value = '1'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.243'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '0.31%'
$ws.Range('F19').NumberFormat = '@'
$ws.Range('F19').Value = '14-1-2023'
$ws.Range('G19').NumberFormat = '@'
$ws.Range('G19').Value = '1'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.3263'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '1.61%'
$ws.Range('F20').NumberFormat = '@'
$ws.Range('F20').Value = '14-1-2023'
$ws.Range('G20').NumberFormat = '@'
$ws.Range('G20').Value = '1'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1318'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '1.53%'
$ws.Range('F21').NumberFormat = '@'
$ws.Range('F21').Value = '14-1-2023'
$ws.Range('G21').NumberFormat = '@'
$ws.Range('G21').Value = '1'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.091'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '0.39%'
$ws.Range('F22').NumberFormat = '@'
$ws.Range('F22').Value = '14-1-2023'
$ws.Range('G22').NumberFormat = '@'
$ws.Range('G22').Value = '1'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.1569'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-0.07%'
$ws.Range('F23').NumberFormat = '@'
$ws.Range('F23').Value = '14-1-2023'
$ws.Range('G23').NumberFormat = '@'
$ws.Range('G23').Value = '1'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.04444'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-1.16%'
$ws.Range('F24').NumberFormat = '@'
$ws.Range('F24').Value = '14-1-2023'
$ws.Range('G24').NumberFormat = '@'
$ws.Range('G24').Value = '1'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.001194'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '0.18%'
$ws.Range('F25').NumberFormat = '@'
$ws.Range('F25').Value = '14-1-2023'
$ws.Range('G25').NumberFormat = '@'
$ws.Range('G25').Value = '1'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.004479'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '8.89%'
$ws.Range('F26').NumberFormat = '@'
$ws.Range('F26').Value = '14-1-2023'
$ws.Range('G26').NumberFormat = '@'
$ws.Range('G26').Value = '1'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0001153'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '-7.83%'
$ws.Range('F27').NumberFormat = '@'
$ws.Range('F27').Value = '14-1-2023'
$ws.Range('G27').NumberFormat = '@'
$ws.Range('G27').Value = '1'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0001672'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '3.13%'
$ws.Range('F28').NumberFormat = '@'
$ws.Range('F28').Value = '14-1-2023'
$ws.Range('G28').NumberFormat = '@'
$ws.Range('G28').Value = '1'
$ws.Range('F29').NumberFormat = '@'
$ws.Range('F29').Value = '14-1-2023'
$ws.Range('G29').NumberFormat = '@'
$ws.Range('G29').Value = '1'
$ws.Range('F30').NumberFormat = '@'
$ws.Range('F30').Value = '14-1-2023'
$ws.Range('G30').NumberFormat = '@'
$ws.Range('G30').Value = '1'
$ws.Range('F31').NumberFormat = '@'
$ws.Range('F31').Value = '14-1-2023'
$ws.Range('G31').NumberFormat = '@'
$ws.Range('G31').Value = '1'
$ws.Range('F32').NumberFormat = '@'
$ws.Range('F32').Value = '14-1-2023'
$ws.Range('G32').NumberFormat = '@'
$ws.Range('G32').Value = '1'
$ws.Range('F33').NumberFormat = '@'
$ws.Range('F33').Value = '14-1-2023'
$ws.Range('G33').NumberFormat = '@'
$ws.Range('G33').Value = '1'
$ws.Range('F34').NumberFormat = '@'
$ws.Range('F34').Value = '14-1-2023'
$ws.Range('G34').NumberFormat = '@'
$ws.Range('G34').Value = '1'
$ws.Range('F35').NumberFormat = '@'
$ws.Range('F35').Value = '14-1-2023'
$ws.Range('G35').NumberFormat = '@'
$ws.Range('G35').Value = '1'
$ws.Range('F36').NumberFormat = '@'
$ws.Range('F36').Value = '14-1-2023'
$ws.Range('G36').NumberFormat = '@'
$ws.Range('G36').Value = '1'
$ws.Range('F37').NumberFormat = '@'
$ws.Range('F37').Value = '14-1-2023'
$ws.Range('G37').NumberFormat = '@'
$ws.Range('G37').Value = '1'
$ws.Range('F38').NumberFormat = '@'
$ws.Range('F38').Value = '14-1-2023'
$ws.Range('G38').NumberFormat = '@'
$ws.Range('G38').Value = '1'
$ws.Range('F39').NumberFormat = '@'
$ws.Range('F39').Value = '14-1-2023'
$ws.Range('G39').NumberFormat = '@'
$ws.Range('G39').Value = '1'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04391'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '4.89%'
$ws.Range('F40').NumberFormat = '@'
$ws.Range('F40').Value = '14-1-2023'
$ws.Range('G40').NumberFormat = '@'
$ws.Range('G40').Value = '1'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006656'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-1.31%'
$ws.Range('F41').NumberFormat = '@'
$ws.Range('F41').Value = '14-1-2023'
$ws.Range('G41').NumberFormat = '@'
$ws.Range('G41').Value = '1'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1344'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '8.37%'
$ws.Range('F42').NumberFormat = '@'
$ws.Range('F42').Value = '14-1-2023'
$ws.Range('G42').NumberFormat = '@'
$ws.Range('G42').Value = '1'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002147'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '8.36%'
$ws.Range('F43').NumberFormat = '@'
$ws.Range('F43').Value = '14-1-2023'
$ws.Range('G43').NumberFormat = '@'
$ws.Range('G43').Value = '1'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.01305'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '4.66%'
$ws.Range('F44').NumberFormat = '@'
$ws.Range('F44').Value = '14-1-2023'
$ws.Range('G44').NumberFormat = '@'
$ws.Range('G44').Value = '1'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005733'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '1.57%'
$ws.Range('F45').NumberFormat = '@'
$ws.Range('F45').Value = '14-1-2023'
$ws.Range('G45').NumberFormat = '@'
$ws.Range('G45').Value = '1'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.885'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-4.17%'
$ws.Range('F46').NumberFormat = '@'
$ws.Range('F46').Value = '14-1-2023'
$ws.Range('G46').NumberFormat = '@'
$ws.Range('G46').Value = '1'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.01281'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-2.20%'
$ws.Range('F47').NumberFormat = '@'
$ws.Range('F47').Value = '14-1-2023'
$ws.Range('G47').NumberFormat = '@'
$ws.Range('G47').Value = '1'
$ws.Range('F48').NumberFormat = '@'
$ws.Range('F48').Value = '14-1-2023'
$ws.Range('G48').NumberFormat = '@'
$ws.Range('G48').Value = '1'
$ws.Range('F49').NumberFormat = '@'
$ws.Range('F49').Value = '14-1-2023'
$ws.Range('G49').NumberFormat = '@'
$ws.Range('G49').Value = '1'
$ws.Range('F50').NumberFormat = '@'
$ws.Range('F50').Value = '14-1-2023'
$ws.Range('G50').NumberFormat = '@'
$ws.Range('G50').Value = '1'
$ws.Range('F51').NumberFormat = '@'
$ws.Range('F51').Value = '14-1-2023'
$ws.Range('G51').NumberFormat = '@'
$ws.Range('G51').Value = '1'
